# Daily attendance processing - swap "Recorded By" name order from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System" wherever it
# appears in the "Recorded By" column (column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$col = $ws.Columns(7)

$first = $col.Find($target)
if ($first) {
    $firstAddress = $first.Address()
    $cell = $first
    $guard = 0
    do {
        $cell.Value2 = $replacement
        $guard = $guard + 1
        $cell = $col.FindNext($cell)
    } while ($cell -and $cell.Address() -ne $firstAddress -and $guard -lt 1000)
}
